$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 49/50 swap: Algorand <-> EnergySwap ---
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

# --- D column: text-safe values (contain multiple dots, never parsed as numbers) ---
$ws.Cells.Item(2, 4).Value = '27.251.56'
$ws.Cells.Item(3, 4).Value = '1.656.30'
$ws.Cells.Item(12, 4).Value = '1.883.39'
$ws.Cells.Item(13, 4).Value = '1.655.51'
$ws.Cells.Item(17, 4).Value = '27.199.76'
$ws.Cells.Item(18, 4).Value = '0.0₃0739'
$ws.Cells.Item(35, 4).Value = '1.266.85'
$ws.Cells.Item(43, 4).Value = '1.792.98'

# --- D column: numeric-looking values that must remain plain text ---
# (preserve original cell style: force Text format, set value, restore style)
$dNumeric = [ordered]@{
    5 = '220.24'
    6 = '0.503'
    8 = '0.255'
    10 = '19.66'
    11 = '0.0848'
    14 = '4.21'
    16 = '66.12'
    19 = '221.81'
    21 = '6.77'
    22 = '4.44'
    23 = '2.44'
    24 = '9.29'
    25 = '147.54'
    27 = '7.37'
    29 = '15.93'
    33 = '3.02'
    38 = '0.540'
    39 = '0.829'
    41 = '0.809'
    42 = '5.39'
    44 = '62.05'
    45 = '92.70'
    46 = '2.09'
    47 = '1.62'
    49 = '7.66'
    50 = '0.0978'
    51 = '0.407'
}
foreach ($row in $dNumeric.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $dNumeric[$row]
    $cell.Style = $origStyle
}

# --- E column: percentage-change text values (safe to set directly) ---
$eChanges = [ordered]@{
    2 = '  +1.52%  '
    3 = '  +0.60%  '
    4 = '  -0.82%  '
    5 = '  +1.59%  '
    6 = '  +0.00%  '
    7 = '  -0.81%  '
    8 = '  +0.88%  '
    9 = '  -0.20%  '
    10 = '  +2.26%  '
    11 = '  +0.42%  '
    12 = '  +0.31%  '
    13 = '  +0.57%  '
    14 = '  +1.04%  '
    15 = '  +0.54%  '
    16 = '  +2.08%  '
    17 = '  +1.30%  '
    18 = '  +0.39%  '
    19 = '  +3.42%  '
    20 = '  -0.77%  '
    21 = '  +8.11%  '
    22 = '  +0.93%  '
    23 = '  -1.84%  '
    24 = '  -0.71%  '
    25 = '  +0.33%  '
    26 = '  -0.79%  '
    27 = '  +2.28%  '
    28 = '  +0.25%  '
    29 = '  +1.51%  '
    30 = '  +1.52%  '
    31 = '  +0.84%  '
    32 = '  +0.69%  '
    33 = '  -0.10%  '
    34 = '  +2.83%  '
    35 = '  -2.38%  '
    36 = '  -0.12%  '
    37 = '  -1.40%  '
    38 = '  +0.89%  '
    39 = '  +0.62%  '
    40 = '  -0.68%  '
    41 = '  +0.03%  '
    42 = '  +1.01%  '
    44 = '  +0.72%  '
    45 = '  +0.75%  '
    46 = '  -6.90%  '
    47 = '  +0.32%  '
    48 = '  -0.70%  '
    49 = '  -0.28%  '
    50 = '  +0.87%  '
    51 = '  -0.16%  '
}
foreach ($row in $eChanges.Keys) {
    $ws.Cells.Item($row, 5).Value = $eChanges[$row]
}
